$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at row 26, shifting existing rows 26-98 down to 35-107
$ws.Range("A26:A34").EntireRow.Insert()

# Force text-formatted columns (B: date string, C: id string) so the
# engine stores them as text, matching the surrounding rows, then restore
# the default "Normal" style so no stray formatting is left behind.
$ws.Range("B26:C34").NumberFormat = "@"

$ws.Range("A26").Value = 1574035200
$ws.Range("B26").Value = "2019-11-18"
$ws.Range("C26").Value = "5293"
$ws.Range("D26").Value = "AME"
$ws.Range("E26").Value = 1.83
$ws.Range("F26").Value = 1.87
$ws.Range("G26").Value = 1.81
$ws.Range("H26").Value = 1.83
$ws.Range("I26").Value = 5107800

$ws.Range("A27").Value = 1574121600
$ws.Range("B27").Value = "2019-11-19"
$ws.Range("C27").Value = "5293"
$ws.Range("D27").Value = "AME"
$ws.Range("E27").Value = 1.82
$ws.Range("F27").Value = 1.85
$ws.Range("G27").Value = 1.82
$ws.Range("H27").Value = 1.83
$ws.Range("I27").Value = 2032800

$ws.Range("A28").Value = 1574208000
$ws.Range("B28").Value = "2019-11-20"
$ws.Range("C28").Value = "5293"
$ws.Range("D28").Value = "AME"
$ws.Range("E28").Value = 1.82
$ws.Range("F28").Value = 1.87
$ws.Range("G28").Value = 1.81
$ws.Range("H28").Value = 1.85
$ws.Range("I28").Value = 3644600

$ws.Range("A29").Value = 1574294400
$ws.Range("B29").Value = "2019-11-21"
$ws.Range("C29").Value = "5293"
$ws.Range("D29").Value = "AME"
$ws.Range("E29").Value = 1.85
$ws.Range("F29").Value = 1.9
$ws.Range("G29").Value = 1.85
$ws.Range("H29").Value = 1.9
$ws.Range("I29").Value = 6461300

$ws.Range("A30").Value = 1574380800
$ws.Range("B30").Value = "2019-11-22"
$ws.Range("C30").Value = "5293"
$ws.Range("D30").Value = "AME"
$ws.Range("E30").Value = 1.9
$ws.Range("F30").Value = 1.91
$ws.Range("G30").Value = 1.87
$ws.Range("H30").Value = 1.88
$ws.Range("I30").Value = 2616400

$ws.Range("A31").Value = 1574640000
$ws.Range("B31").Value = "2019-11-25"
$ws.Range("C31").Value = "5293"
$ws.Range("D31").Value = "AME"
$ws.Range("E31").Value = 1.89
$ws.Range("F31").Value = 1.95
$ws.Range("G31").Value = 1.88
$ws.Range("H31").Value = 1.92
$ws.Range("I31").Value = 4421000

$ws.Range("A32").Value = 1574726400
$ws.Range("B32").Value = "2019-11-26"
$ws.Range("C32").Value = "5293"
$ws.Range("D32").Value = "AME"
$ws.Range("E32").Value = 1.93
$ws.Range("F32").Value = 1.96
$ws.Range("G32").Value = 1.91
$ws.Range("H32").Value = 1.94
$ws.Range("I32").Value = 2286900

$ws.Range("A33").Value = 1574812800
$ws.Range("B33").Value = "2019-11-27"
$ws.Range("C33").Value = "5293"
$ws.Range("D33").Value = "AME"
$ws.Range("E33").Value = 1.94
$ws.Range("F33").Value = 1.97
$ws.Range("G33").Value = 1.93
$ws.Range("H33").Value = 1.94
$ws.Range("I33").Value = 5141500

$ws.Range("A34").Value = 1574899200
$ws.Range("B34").Value = "2019-11-28"
$ws.Range("C34").Value = "5293"
$ws.Range("D34").Value = "AME"
$ws.Range("E34").Value = 1.95
$ws.Range("F34").Value = 1.96
$ws.Range("G34").Value = 1.92
$ws.Range("H34").Value = 1.92
$ws.Range("I34").Value = 2619700

$ws.Range("B26:C34").Style = "Normal"
